# Ahila final project updates
# Adds two new worksheets to the end of the workbook:
#   - UserDeleteByUserId         (valid userId test data)
#   - UserDeleteByInvalidUserId  (invalid userId test data, with a hyperlinked value)

$wb = $excel.ActiveWorkbook

# Anchor on the current last sheet (BatchDeleteByValidId) so the new sheets are appended at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- New sheet: UserDeleteByUserId ---
$userDeleteByUserId = $wb.Worksheets.Add($null, $lastSheet)
$userDeleteByUserId.Name = "UserDeleteByUserId"
$userDeleteByUserId.Range("A1").Value = "userId"
$userDeleteByUserId.Range("A2").Value = "U123454"
$userDeleteByUserId.Activate() | Out-Null
$userDeleteByUserId.Range("A2").Select() | Out-Null

# --- New sheet: UserDeleteByInvalidUserId ---
$userDeleteByInvalidUserId = $wb.Worksheets.Add($null, $userDeleteByUserId)
$userDeleteByInvalidUserId.Name = "UserDeleteByInvalidUserId"
$userDeleteByInvalidUserId.Range("A1").Value = "userId"
$userDeleteByInvalidUserId.Range("A2").Value = "U132@"

# The invalid userId value is rendered as a hyperlink (Excel auto-creates the
# built-in "Hyperlink" cell style/font the first time a hyperlink is added).
$userDeleteByInvalidUserId.Hyperlinks.Add($userDeleteByInvalidUserId.Range("A2"), "https://example.com/users/U132@") | Out-Null

# This becomes the active sheet/tab in the saved workbook.
$userDeleteByInvalidUserId.Activate() | Out-Null
$userDeleteByInvalidUserId.Range("A2").Select() | Out-Null
